$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text before writing numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.918.21'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '1.670.89'
$ws.Range("E3").Value = '  +0.99%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '214.88'
$ws.Range("E5").Value = '  -0.04%  '

$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.250'
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.0620'
$ws.Range("E9").Value = '  +0.81%  '

$ws.Range("D10").Value = '20.26'
$ws.Range("E10").Value = '  +0.12%  '

$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +1.77%  '

$ws.Range("D12").Value = '1.905.81'
$ws.Range("E12").Value = '  +0.94%  '

$ws.Range("D13").Value = '1.655.31'
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("E14").Value = '  +0.04%  '

$ws.Range("E15").Value = '  +1.42%  '

$ws.Range("D16").Value = '65.49'
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("D17").Value = '26.911.62'
$ws.Range("E17").Value = '  -0.14%  '

$ws.Range("D18").Value = '8.03'
$ws.Range("E18").Value = '  +3.99%  '

$ws.Range("D19").Value = '233.26'
$ws.Range("E19").Value = '  -1.31%  '

$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("E21").Value = '  -0.07%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").Value = '9.14'
$ws.Range("E23").Value = '  -1.64%  '

$ws.Range("E24").Value = '  -2.11%  '

$ws.Range("D25").Value = '146.00'
$ws.Range("E25").Value = '  +0.45%  '

$ws.Range("D26").Value = '7.12'
$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("D27").Value = '15.95'
$ws.Range("E27").Value = '  +0.96%  '

$ws.Range("E29").Value = '  -1.86%  '

$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  +0.23%  '

$ws.Range("E31").Value = '  +0.04%  '

$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").Value = '1.457.58'
$ws.Range("E33").Value = '  -5.59%  '

$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("E35").Value = '  +1.75%  '

$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").Value = '0.581'
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("E38").Value = '  +0.87%  '

$ws.Range("E39").Value = '  +0.39%  '

$ws.Range("E40").Value = '  +13.00%  '

$ws.Range("D41").Value = '5.74'
$ws.Range("E41").Value = '  -4.41%  '

$ws.Range("E42").Value = '  -0.06%  '

$ws.Range("E43").Value = '  +2.12%  '

$ws.Range("D44").Value = '66.23'
$ws.Range("E44").Value = '  +0.80%  '

$ws.Range("D45").Value = '1.811.62'
$ws.Range("E45").Value = '  +0.92%  '

$ws.Range("E46").Value = '  +0.71%  '

$ws.Range("D47").Value = '90.73'
$ws.Range("E47").Value = '  +0.46%  '

$ws.Range("E48").Value = '  +1.14%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.101'
$ws.Range("E49").Value = '  +2.69%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0507'
$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.62'
$ws.Range("E51").Value = '  -0.02%  '

# Reset style on column D back to default (Normal) while keeping text type
$ws.Range("D2:D51").Style = "Normal"
